# feat: add 2022-Q3 data
#
# Insert a new "2022-Q3" worksheet right after "总计" (pushing 2022-Q1 and
# 2021-Q4 one position to the right), and update the "总计" summary sheet
# so its top data row now reflects the new 2022-Q3 totals, with the old
# 2022-Q1 / 2021-Q4 rows shifted down by one.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Item("2022-Q1")

# --- 1. Insert the new "2022-Q3" sheet right after "总计" --------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Header row (same headers/layout as the other quarterly sheets)
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Fund holding rows
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "519615"
$q3.Range("C2").Value = "银河君尚灵活配置混合I"
$q3.Range("D2").Value = "3.59"
$q3.Range("E2").Value = "35.36"
$q3.Range("F2").Value = "0.72"
$q3.Range("G2").Value = "0.0258"
$q3.Range("H2").Value = 9

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "519613"
$q3.Range("C3").Value = "银河君尚灵活配置混合A"
$q3.Range("D3").Value = "2.10"
$q3.Range("E3").Value = "35.36"
$q3.Range("F3").Value = "0.72"
$q3.Range("G3").Value = "0.0151"
$q3.Range("H3").Value = 9

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "519614"
$q3.Range("C4").Value = "银河君尚灵活配置混合C"
$q3.Range("D4").Value = "0.17"
$q3.Range("E4").Value = "35.36"
$q3.Range("F4").Value = "0.72"
$q3.Range("G4").Value = "0.0012"
$q3.Range("H4").Value = 9

# --- 2. Update the "总计" sheet: shift the existing rows down and add --
#        the new 2022-Q3 totals on top ------------------------------------

# Make room: push the old row 2 (2022-Q1) / row 3 (2021-Q4) down to rows 3/4
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q4"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 0.01

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q1"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.07000000000000001

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.04
